$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the figures in columns D and E for rows 2-7
$ws.Range("D2").Value = 447981
$ws.Range("E2").Value = 300

$ws.Range("D3").Value = -297573
$ws.Range("E3").Value = 6527

$ws.Range("D4").Value = 43484
$ws.Range("E4").Value = 198

$ws.Range("D5").Value = 50400

$ws.Range("D6").Value = 10492

$ws.Range("D7").Value = 26130

# Move / record the active selection as it was left in the saved file
$ws.Range("F15").Select() | Out-Null
